$d = $word.ActiveDocument

$d.Content.Find.Execute("Coronavirus Scale : 11.67 %", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Coronavirus Scale : 91.53 %", 2)

$d.Content.Find.Execute("Automated Tests : Passed Successfully", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Automated Tests : Failed", 2)

$d.Content.Find.Execute("Manual Tests : Passed with Considerations", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Manual Tests : Failed, Found COVID Positive", 2)

$d.Content.Find.Execute("Comments : Fit for Travel", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Comments : Unfit for Travel", 2)

$d.Content.Find.Execute("Generated : 22.06.2020 5:30GMT", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Generated : 19.04.2020 3:30GMT", 2)
